# Benfords_law_covid_cases: add a "last digit" percentage column alongside
# the existing "first digit" percentage column, per commit
# "Added last digit percentage column and snapshots".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing column C ("Covid_dataset_percentage") already holds the
# covid-dataset FIRST-digit percentages; re-label it to be explicit now
# that a last-digit counterpart is being added. Column D
# ("Benfords_law_percentage") and its values are unchanged.
$ws.Range("C1").Value = "Covid_dataset_first_digit_percentage"

# New column E: covid-dataset LAST-digit percentages (DIGIT 1..9).
$ws.Range("E1").Value = "Covid_dataset_last_digit_percentage"
$ws.Range("E2").Value = 13.4093
$ws.Range("E3").Value = 12.118
$ws.Range("E4").Value = 11.6833
$ws.Range("E5").Value = 10.9939
$ws.Range("E6").Value = 10.724
$ws.Range("E7").Value = 10.507
$ws.Range("E8").Value = 10.2999
$ws.Range("E9").Value = 10.1943
$ws.Range("E10").Value = 10.0705

# Resize the columns to fit their (now longer) headers/content.
$ws.Columns("A:E").AutoFit()

# Match the author's saved cursor position.
$ws.Range("F12").Select() | Out-Null
